# [UPDATE] Perubahan format xlsx
#
# 1. Month headers in row 1 (C1:H1) change from "Jan  2018".."Jun  2018"
#    to the short numeric form "1-2018".."6-2018", and are formatted as
#    Text (numFmtId 49 / "@") so the leading digit isn't reinterpreted.
# 2. The little "Nama Bulan" lookup table (B24 label + C24:C35 month
#    abbreviations) is cleared out - it's no longer used.
# 3. The active selection moves from C28 to L22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rewrite the month/year headers and mark them as Text ---
$months = @("1-2018", "2-2018", "3-2018", "4-2018", "5-2018", "6-2018")
$cols = @("C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.Value = $months[$i]
    $cell.NumberFormat = "@"
}

# --- 2. Clear the now-unused "Nama Bulan" helper table ---
$ws.Range("B24").ClearContents()
$ws.Range("C24:C35").ClearContents()

# --- 3. Move the selection ---
$ws.Range("L22").Select()
